$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at row 648, pushing existing data (old rows 648-658) down to 654-664
$ws.Rows("648:653").Insert()

# Populate the newly inserted rows with the new weekly data
# Row 648
$ws.Range("A648").Value = 3
$ws.Range("B648").Value = 'Femacal de La Calera'
$ws.Range("C648").Value = 'Coquimbo'
$ws.Range("D648").Value = 44448
$ws.Range("E648").Value = 5
$ws.Range("F648").Value = 'Fruta'
$ws.Range("G648").Value = 100104
$ws.Range("H648").Value = 'Frutos de pepita'
$ws.Range("I648").Value = 100104005
$ws.Range("J648").Value = 'Pera'
$ws.Range("K648").Value = 'Packham''s Triumph'
$ws.Range("L648").Value = 'Especial'
$ws.Range("M648").Value = 130
$ws.Range("N648").Value = 12000
$ws.Range("O648").Value = 12000
$ws.Range("P648").Value = 12000
$ws.Range("Q648").Value = '$/caja 18 kilos empedrada'
$ws.Range("R648").Value = 'Provincia de Colchagua'
$ws.Range("S648").Value = 667
$ws.Range("T648").Value = 18

# Row 649
$ws.Range("A649").Value = 3
$ws.Range("B649").Value = 'Femacal de La Calera'
$ws.Range("C649").Value = 'Coquimbo'
$ws.Range("D649").Value = 44448
$ws.Range("E649").Value = 5
$ws.Range("F649").Value = 'Fruta'
$ws.Range("G649").Value = 100104
$ws.Range("H649").Value = 'Frutos de pepita'
$ws.Range("I649").Value = 100104005
$ws.Range("J649").Value = 'Pera'
$ws.Range("K649").Value = 'Packham''s Triumph'
$ws.Range("L649").Value = 'Extra (doble especial)'
$ws.Range("M649").Value = 120
$ws.Range("N649").Value = 13000
$ws.Range("O649").Value = 13000
$ws.Range("P649").Value = 13000
$ws.Range("Q649").Value = '$/caja 18 kilos empedrada'
$ws.Range("R649").Value = 'Provincia de Colchagua'
$ws.Range("S649").Value = 722
$ws.Range("T649").Value = 18

# Row 650
$ws.Range("A650").Value = 3
$ws.Range("B650").Value = 'Femacal de La Calera'
$ws.Range("C650").Value = 'Coquimbo'
$ws.Range("D650").Value = 44448
$ws.Range("E650").Value = 5
$ws.Range("F650").Value = 'Fruta'
$ws.Range("G650").Value = 100104
$ws.Range("H650").Value = 'Frutos de pepita'
$ws.Range("I650").Value = 100104005
$ws.Range("J650").Value = 'Pera'
$ws.Range("K650").Value = 'Packham''s Triumph'
$ws.Range("L650").Value = 'Primera'
$ws.Range("M650").Value = 100
$ws.Range("N650").Value = 11000
$ws.Range("O650").Value = 11000
$ws.Range("P650").Value = 11000
$ws.Range("Q650").Value = '$/caja 18 kilos empedrada'
$ws.Range("R650").Value = 'Provincia de Colchagua'
$ws.Range("S650").Value = 611
$ws.Range("T650").Value = 18

# Row 651
$ws.Range("A651").Value = 3
$ws.Range("B651").Value = 'Femacal de La Calera'
$ws.Range("C651").Value = 'Coquimbo'
$ws.Range("D651").Value = 44448
$ws.Range("E651").Value = 5
$ws.Range("F651").Value = 'Fruta'
$ws.Range("G651").Value = 100104
$ws.Range("H651").Value = 'Frutos de pepita'
$ws.Range("I651").Value = 100104005
$ws.Range("J651").Value = 'Pera'
$ws.Range("K651").Value = 'Packham''s Triumph'
$ws.Range("L651").Value = 'Segunda'
$ws.Range("M651").Value = 115
$ws.Range("N651").Value = 10000
$ws.Range("O651").Value = 10000
$ws.Range("P651").Value = 10000
$ws.Range("Q651").Value = '$/caja 18 kilos empedrada'
$ws.Range("R651").Value = 'Provincia de Colchagua'
$ws.Range("S651").Value = 556
$ws.Range("T651").Value = 18

# Row 652
$ws.Range("A652").Value = 3
$ws.Range("B652").Value = 'Femacal de La Calera'
$ws.Range("C652").Value = 'Coquimbo'
$ws.Range("D652").Value = 44448
$ws.Range("E652").Value = 5
$ws.Range("F652").Value = 'Fruta'
$ws.Range("G652").Value = 100104
$ws.Range("H652").Value = 'Frutos de pepita'
$ws.Range("I652").Value = 100104005
$ws.Range("J652").Value = 'Pera'
$ws.Range("K652").Value = 'Winter Nelis'
$ws.Range("L652").Value = 'Primera'
$ws.Range("M652").Value = 75
$ws.Range("N652").Value = 11000
$ws.Range("O652").Value = 11000
$ws.Range("P652").Value = 11000
$ws.Range("Q652").Value = '$/caja 18 kilos empedrada'
$ws.Range("R652").Value = 'Provincia de Colchagua'
$ws.Range("S652").Value = 611
$ws.Range("T652").Value = 18

# Row 653
$ws.Range("A653").Value = 3
$ws.Range("B653").Value = 'Femacal de La Calera'
$ws.Range("C653").Value = 'Coquimbo'
$ws.Range("D653").Value = 44448
$ws.Range("E653").Value = 5
$ws.Range("F653").Value = 'Fruta'
$ws.Range("G653").Value = 100104
$ws.Range("H653").Value = 'Frutos de pepita'
$ws.Range("I653").Value = 100104005
$ws.Range("J653").Value = 'Pera'
$ws.Range("K653").Value = 'Winter Nelis'
$ws.Range("L653").Value = 'Segunda'
$ws.Range("M653").Value = 70
$ws.Range("N653").Value = 10000
$ws.Range("O653").Value = 10000
$ws.Range("P653").Value = 10000
$ws.Range("Q653").Value = '$/caja 18 kilos empedrada'
$ws.Range("R653").Value = 'Provincia de Colchagua'
$ws.Range("S653").Value = 556
$ws.Range("T653").Value = 18
